$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.872.08"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.26%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.111.47"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.72%  "

$ws.Range("E4").Value = "  +0.04%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "575.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.54%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.30"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.94%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.106.39"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.72%  "

$ws.Range("E9").Value = "  -0.44%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.41"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.54%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.152"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.41%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.477"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.99%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000245"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.89%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "37.08"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.30%  "

$ws.Range("E15").Value = "  -1.00%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.627.16"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.58%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.926.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.26%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.09"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.09%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.111.57"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.65%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.16%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "476.00"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +3.27%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.710"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.74"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.69%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "83.73"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.97%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "13.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.53%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.28"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.83%  "

$ws.Range("E27").Value = "  -0.05%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.91"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.33%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.98"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.37%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.41"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.53%  "

$ws.Range("E31").Value = "  -0.05%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.55"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.34%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0₃0963"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -5.05%  "

$ws.Range("E34").Value = "  -1.35%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.04%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.84"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.02%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.978"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.55%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "47.82"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.77%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.08"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +2.65%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "50.04"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.48%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.309"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.46%  "

$ws.Range("E42").Value = "  -0.16%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.59"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.42%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.795.66"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.52%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0354"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.87%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "377.79"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.43%  "

$ws.Range("E47").Value = "  -9.99%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "136.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.76%  "

$ws.Range("E49").Value = "  -0.03%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "24.67"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.16%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.19"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.90%  "

